$d = $word.ActiveDocument

$pairs = @(
    @("643÷4=", "861÷9="),
    @("382÷7=", "645÷9="),
    @("774÷8=", "755÷6="),
    @("948÷5=", "143÷3="),
    @("951÷7=", "609÷9="),
    @("289÷2=", "724÷8="),
    @("788÷2=", "907÷3="),
    @("283÷2=", "571÷9="),
    @("466÷9=", "235÷6="),
    @("801÷2=", "684÷6="),
    @("446÷8=", "449÷4="),
    @("360÷3=", "598÷5="),
    @("457÷5=", "758÷7="),
    @("684÷5=", "109÷4="),
    @("817÷3=", "361÷6="),
    @("899÷8=", "276÷4="),
    @("666÷6=", "986÷6="),
    @("900÷8=", "696÷8="),
    @("261÷2=", "623÷2="),
    @("469÷3=", "254÷7="),
    @("249÷3=", "740÷9="),
    @("193÷5=", "500÷4="),
    @("483÷7=", "950÷8="),
    @("802÷2=", "838÷6="),
    @("806÷4=", "865÷7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
